$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-03 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-04 Sunday", 2) | Out-Null
$d.Content.Find.Execute("460÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "789÷7=", 2) | Out-Null
$d.Content.Find.Execute("402÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "909÷2=", 2) | Out-Null
$d.Content.Find.Execute("668÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "191÷9=", 2) | Out-Null
$d.Content.Find.Execute("869÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "249÷3=", 2) | Out-Null
$d.Content.Find.Execute("605÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "753÷7=", 2) | Out-Null
$d.Content.Find.Execute("975÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "215÷8=", 2) | Out-Null
$d.Content.Find.Execute("462÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "331÷6=", 2) | Out-Null
$d.Content.Find.Execute("825÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "403÷5=", 2) | Out-Null
$d.Content.Find.Execute("729÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "247÷6=", 2) | Out-Null
$d.Content.Find.Execute("321÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "677÷4=", 2) | Out-Null
$d.Content.Find.Execute("147÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "258÷8=", 2) | Out-Null
$d.Content.Find.Execute("999÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷3=", 2) | Out-Null
$d.Content.Find.Execute("889÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "937÷2=", 2) | Out-Null
$d.Content.Find.Execute("718÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "885÷5=", 2) | Out-Null
$d.Content.Find.Execute("231÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "542÷3=", 2) | Out-Null
$d.Content.Find.Execute("737÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "610÷3=", 2) | Out-Null
$d.Content.Find.Execute("264÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "903÷2=", 2) | Out-Null
$d.Content.Find.Execute("672÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷9=", 2) | Out-Null
$d.Content.Find.Execute("294÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "627÷2=", 2) | Out-Null
$d.Content.Find.Execute("862÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "709÷8=", 2) | Out-Null
$d.Content.Find.Execute("678÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "942÷6=", 2) | Out-Null
$d.Content.Find.Execute("639÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "489÷4=", 2) | Out-Null
$d.Content.Find.Execute("642÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "778÷9=", 2) | Out-Null
$d.Content.Find.Execute("773÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "254÷8=", 2) | Out-Null
$d.Content.Find.Execute("345÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "234÷4=", 2) | Out-Null
